$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = 0.99999999381659688
$ws.Range("A2").Value = 0.99798790574396001
$ws.Range("A3").Value = 0.99582741435549349
$ws.Range("A4").Value = 1.0001165426063676
$ws.Range("A5").Value = 0.99130890133367267
$ws.Range("A6").Value = 0.97052642641708609
$ws.Range("A7").Value = 0.96930913890157899
$ws.Range("A8").Value = 0.96768662457505439
$ws.Range("A9").Value = 0.96915930392030958
$ws.Range("A10").Value = 0.96064686061925608
$ws.Range("A11").Value = 0.95956916201136477
$ws.Range("A12").Value = 0.95777511981127739
$ws.Range("A13").Value = 0.94956451553613563
$ws.Range("A14").Value = 0.94539786879675036
$ws.Range("A15").Value = 0.94280672461487325
$ws.Range("A16").Value = 0.94030046813251544
$ws.Range("A17").Value = 0.93659276215077403
$ws.Range("A18").Value = 0.93548385576702842
$ws.Range("A19").Value = 0.99474986018500489
$ws.Range("A20").Value = 0.96834168784077046
$ws.Range("A21").Value = 0.96694317389374973
$ws.Range("A22").Value = 0.9656786592180957
$ws.Range("A23").Value = 0.98879340286076745
$ws.Range("A24").Value = 0.97577332987428689
$ws.Range("A25").Value = 0.96931646806051475
$ws.Range("A26").Value = 0.9673340353272033
$ws.Range("A27").Value = 0.96313156199979533
$ws.Range("A28").Value = 0.94533370122376215
$ws.Range("A29").Value = 0.93276796250322502
$ws.Range("A30").Value = 0.92727898754929328
$ws.Range("A31").Value = 0.92321628313818127
$ws.Range("A32").Value = 0.92285941045725739
$ws.Range("A33").Value = 0.9223394097840687
